# edit.ps1 -- apply the "Updated cryptos list" data refresh to Sheet1
#
# The worksheet stores every data cell (coin name, link, price, volume) as
# plain text (inline strings) -- including price values such as "229.81"
# that look like ordinary decimal numbers. Excel's COM automation layer
# auto-detects numeric-looking text and converts it to a real number when
# you assign it through Range.Value. To keep those particular cells as text
# (matching the source feed / original formatting, e.g. trailing zeros like
# "14.60"), we briefly switch them to the Text ("@") number format before
# writing the value, then clear the formatting again so the cells end up
# with the same (default/General) style as before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$textForceCells = @(
    "D5",
    "D6",
    "D7",
    "D10",
    "D11",
    "D13",
    "D14",
    "D15",
    "D16",
    "D19",
    "D20",
    "D22",
    "D24",
    "D25",
    "D26",
    "D27",
    "D29",
    "D33",
    "D37",
    "D42",
    "D44",
    "D45",
    "D47"
)
$textForceRange = $ws.Range($textForceCells[0])
foreach ($addr in $textForceCells[1..($textForceCells.Length - 1)]) {
    $textForceRange = $excel.Union($textForceRange, $ws.Range($addr))
}
$textForceRange.NumberFormat = "@"

# Row 2: Bitcoin
$ws.Range("D2").Value = "37.950.43"
$ws.Range("E2").Value = "  +2.14%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.050.70"
$ws.Range("E3").Value = "  +1.25%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5: BNB
$ws.Range("D5").Value = "229.81"
$ws.Range("E5").Value = "  +1.45%  "

# Row 6: XRP
$ws.Range("D6").Value = "0.616"
$ws.Range("E6").Value = "  +2.21%  "

# Row 7: Solana
$ws.Range("D7").Value = "58.53"
$ws.Range("E7").Value = "  +6.14%  "

# Row 8: USDC
$ws.Range("E8").Value = "  -0.01%  "

# Row 9: Cardano
$ws.Range("E9").Value = "  +1.93%  "

# Row 10: Dogecoin
$ws.Range("D10").Value = "0.0813"
$ws.Range("E10").Value = "  +3.41%  "

# Row 11: TRON
$ws.Range("D11").Value = "0.103"
$ws.Range("E11").Value = "  +2.37%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.353.83"
$ws.Range("E12").Value = "  +1.54%  "

# Row 13: Chainlink
$ws.Range("D13").Value = "14.60"
$ws.Range("E13").Value = "  +2.26%  "

# Row 14: Avalanche
$ws.Range("D14").Value = "20.85"
$ws.Range("E14").Value = "  +2.81%  "

# Row 15: Polygon (was Polkadot)
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.753"
$ws.Range("E15").Value = "  +1.38%  "

# Row 16: Polkadot (was Polygon)
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "5.29"
$ws.Range("E16").Value = "  +2.04%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "2.054.80"
$ws.Range("E17").Value = "  +1.65%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "37.882.67"
$ws.Range("E18").Value = "  +2.09%  "

# Row 19: Uniswap
$ws.Range("D19").Value = "6.27"
$ws.Range("E19").Value = "  -1.55%  "

# Row 20: Litecoin
$ws.Range("D20").Value = "69.79"
$ws.Range("E20").Value = "  +1.33%  "

# Row 21: ShibaInu
$ws.Range("D21").Value = "0.0₃0836"
$ws.Range("E21").Value = "  +2.36%  "

# Row 22: BitcoinCash
$ws.Range("D22").Value = "224.23"
$ws.Range("E22").Value = "  +0.33%  "

# Row 23: Dai
$ws.Range("E23").Value = "  -0.02%  "

# Row 24: Toncoin
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  -0.32%  "

# Row 25: PancakeSwap
$ws.Range("D25").Value = "2.24"
$ws.Range("E25").Value = "  +2.75%  "

# Row 26: Cosmos
$ws.Range("D26").Value = "9.33"
$ws.Range("E26").Value = "  +0.74%  "

# Row 27: Monero
$ws.Range("D27").Value = "166.36"
$ws.Range("E27").Value = "  -0.01%  "

# Row 28: Kaspa
$ws.Range("E28").Value = "  +4.65%  "

# Row 29: EthereumClassic
$ws.Range("D29").Value = "18.98"
$ws.Range("E29").Value = "  +1.37%  "

# Row 30: ImmutableX
$ws.Range("E30").Value = "  +1.04%  "

# Row 31: Stellar
$ws.Range("E31").Value = "  +1.39%  "

# Row 32: Filecoin
$ws.Range("E32").Value = "  +0.21%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").Value = "4.57"
$ws.Range("E33").Value = "  +2.48%  "

# Row 34: Hedera
$ws.Range("E34").Value = "  -0.20%  "

# Row 35: WEMIXToken
$ws.Range("E35").Value = "  +8.99%  "

# Row 36: LidoDAOToken
$ws.Range("E36").Value = "  -1.24%  "

# Row 37: THORChain
$ws.Range("D37").Value = "6.10"
$ws.Range("E37").Value = "  +9.16%  "

# Row 38: RenderToken
$ws.Range("E38").Value = "  +6.17%  "

# Row 39: BinanceUSD
$ws.Range("E39").Value = "  -0.08%  "

# Row 40: VeChain
$ws.Range("E40").Value = "  +1.00%  "

# Row 41: Maker
$ws.Range("D41").Value = "1.479.99"
$ws.Range("E41").Value = "  +0.59%  "

# Row 42: Aave
$ws.Range("D42").Value = "97.17"
$ws.Range("E42").Value = "  +1.34%  "

# Row 43: HuobiToken
$ws.Range("E43").Value = "  +3.81%  "

# Row 44: InjectiveProtocol
$ws.Range("D44").Value = "16.53"
$ws.Range("E44").Value = "  +0.83%  "

# Row 45: Cronos
$ws.Range("D45").Value = "0.0924"
$ws.Range("E45").Value = "  +1.27%  "

# Row 46: TrustWalletToken
$ws.Range("E46").Value = "  -1.15%  "

# Row 47: FTXToken
$ws.Range("D47").Value = "4.14"
$ws.Range("E47").Value = "  +16.45%  "

# Row 48: ARBITRUM
$ws.Range("E48").Value = "  +0.52%  "

# Row 49: MXToken
$ws.Range("E49").Value = "  +1.35%  "

# Row 50: FraxShare
$ws.Range("E50").Value = "  -2.74%  "

# Row 51: RocketPoolETH
$ws.Range("D51").Value = "2.243.32"
$ws.Range("E51").Value = "  +1.82%  "

# Restore the default (General) style on the cells we temporarily switched
# to Text format, so only the cell contents change -- not their formatting.
$textForceRange.ClearFormats()

